$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 331.83334
$ws.Range("I33").Value = 331.83334
$ws.Range("K33").Value = 331.83334
$ws.Range("M33").Value = -102.83334
$ws.Range("H40").Value = 125002350
$ws.Range("J40").Value = 250001950
$ws.Range("L40").Value = 250001950
$ws.Range("N40").Value = -250002300
$ws.Range("H43").Value = 2492.4375
$ws.Range("J43").Value = 1534.8572
$ws.Range("L43").Value = 1534.8572
$ws.Range("N43").Value = -1672.8572
$ws.Range("H99").Value = 3195.1
$ws.Range("J99").Value = 5207.8335
$ws.Range("L99").Value = 15623.5005
$ws.Range("N99").Value = -18619.5005
$ws.Range("H137").Value = 2181.05
$ws.Range("I137").Value = 1375.091
$ws.Range("J137").Value = 3166.111
$ws.Range("K137").Value = 4125.272999999999
$ws.Range("L137").Value = 9498.332999999999
$ws.Range("M137").Value = -1575.272999999999
$ws.Range("N137").Value = -14598.333
$ws.Range("H138").Value = 3475.5151
$ws.Range("I138").Value = 2257.3462
$ws.Range("J138").Value = 4267.325
$ws.Range("K138").Value = 6772.0386
$ws.Range("L138").Value = 12801.975
$ws.Range("M138").Value = -1632.0386
$ws.Range("N138").Value = -23081.975
$ws.Range("H141").Value = 27783766
$ws.Range("I141").Value = 38464468
$ws.Range("K141").Value = 115393404
$ws.Range("M141").Value = -115388224

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4499.098
$ws.Range("I32").Value = 4549.695
$ws.Range("J32").Value = 3006.5
$ws.Range("K32").Value = 4549.695
$ws.Range("L32").Value = 3006.5
$ws.Range("M32").Value = -4262.695
$ws.Range("N32").Value = -3580.5
$ws.Range("H74").Value = 2141.8386
$ws.Range("I74").Value = 2168.2068
$ws.Range("K74").Value = 2168.2068
$ws.Range("M74").Value = -1294.2068
$ws.Range("H77").Value = 2141.8386
$ws.Range("I77").Value = 2168.2068
$ws.Range("K77").Value = 10841.034
$ws.Range("M77").Value = -6473.034

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 13911.143
$ws.Range("J20").Value = 3900
$ws.Range("L20").Value = 3900
$ws.Range("N20").Value = -4394
$ws.Range("H99").Value = 2278.6
$ws.Range("I99").Value = 2210.75
$ws.Range("J99").Value = 2550
$ws.Range("K99").Value = 2210.75
$ws.Range("L99").Value = 2550
$ws.Range("M99").Value = -712.75
$ws.Range("N99").Value = -5546
$ws.Range("H134").Value = 2633579.8
$ws.Range("I134").Value = 1820.3334
$ws.Range("J134").Value = 12502678
$ws.Range("K134").Value = 5461.0002
$ws.Range("L134").Value = 37508034
$ws.Range("M134").Value = -2926.0002
$ws.Range("N134").Value = -37513104

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 33335248
$ws.Range("I16").Value = 50001680
$ws.Range("K16").Value = 50001680
$ws.Range("M16").Value = -50001393
$ws.Range("H58").Value = 2863.5557
$ws.Range("I58").Value = 2710.9167
$ws.Range("J58").Value = 3168.8333
$ws.Range("K58").Value = 2710.9167
$ws.Range("L58").Value = 3168.8333
$ws.Range("M58").Value = -2507.9167
$ws.Range("N58").Value = -3574.8333
$ws.Range("H94").Value = 1595.3
$ws.Range("I94").Value = 1905.5
$ws.Range("J94").Value = 1517.75
$ws.Range("K94").Value = 1905.5
$ws.Range("L94").Value = 1517.75
$ws.Range("M94").Value = -1454.5
$ws.Range("N94").Value = -2419.75
$ws.Range("H113").Value = 33335248
$ws.Range("I113").Value = 50001680
$ws.Range("K113").Value = 50001680
$ws.Range("M113").Value = -49999510
$ws.Range("H136").Value = 2863.5557
$ws.Range("I136").Value = 2710.9167
$ws.Range("J136").Value = 3168.8333
$ws.Range("K136").Value = 8132.750100000001
$ws.Range("L136").Value = 9506.499899999999
$ws.Range("M136").Value = -5582.750100000001
$ws.Range("N136").Value = -14606.4999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 12383
$ws.Range("I17").Value = 2600
$ws.Range("J17").Value = 22166
$ws.Range("K17").Value = 7800
$ws.Range("L17").Value = 66498
$ws.Range("M17").Value = -7631
$ws.Range("N17").Value = -66836
$ws.Range("H34").Value = 7966.375
$ws.Range("J34").Value = 12566.4
$ws.Range("L34").Value = 37699.2
$ws.Range("N34").Value = -37867.2
$ws.Range("H39").Value = 20607
$ws.Range("J39").Value = 25583.25
$ws.Range("L39").Value = 76749.75
$ws.Range("N39").Value = -77337.75
$ws.Range("H55").Value = 8163.2
$ws.Range("I55").Value = 2471.5715
$ws.Range("K55").Value = 7414.7145
$ws.Range("M55").Value = -7237.7145
$ws.Range("H61").Value = 6047.1665
$ws.Range("I61").Value = 590
$ws.Range("K61").Value = 1770
$ws.Range("M61").Value = -1555
$ws.Range("H64").Value = 14148.375
$ws.Range("I64").Value = 6880.75
$ws.Range("J64").Value = 21416
$ws.Range("K64").Value = 20642.25
$ws.Range("L64").Value = 64248
$ws.Range("M64").Value = -20372.25
$ws.Range("N64").Value = -64788
$ws.Range("H67").Value = 14148.375
$ws.Range("I67").Value = 6880.75
$ws.Range("J67").Value = 21416
$ws.Range("K67").Value = 20642.25
$ws.Range("L67").Value = 64248
$ws.Range("M67").Value = -19706.25
$ws.Range("N67").Value = -66120
$ws.Range("H106").Value = 13131.3
$ws.Range("J106").Value = 20552.166
$ws.Range("L106").Value = 61656.49800000001
$ws.Range("N106").Value = -63548.49800000001
$ws.Range("H129").Value = 5351
$ws.Range("I129").Value = 3652.3333
$ws.Range("J129").Value = 7262
$ws.Range("K129").Value = 10956.9999
$ws.Range("L129").Value = 21786
$ws.Range("M129").Value = -5956.999899999999
$ws.Range("N129").Value = -31786
$ws.Range("H140").Value = 5366.522
$ws.Range("I140").Value = 3252.2
$ws.Range("K140").Value = 9756.599999999999
$ws.Range("M140").Value = -4576.599999999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3032115.2
$ws.Range("I132").Value = 1767.3871
$ws.Range("K132").Value = 5302.1613
$ws.Range("M132").Value = -2772.1613

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3565.9048
$ws.Range("I132").Value = 2214.8667
$ws.Range("J132").Value = 6943.5
$ws.Range("K132").Value = 6644.6001
$ws.Range("L132").Value = 20830.5
$ws.Range("M132").Value = -4114.6001
$ws.Range("N132").Value = -25890.5
$ws.Range("H136").Value = 1949.2106
$ws.Range("I136").Value = 1862.5714
$ws.Range("K136").Value = 5587.7142
$ws.Range("M136").Value = -3037.7142

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2380.1667
$ws.Range("I81").Value = 2296.2
$ws.Range("K81").Value = 4592.4
$ws.Range("M81").Value = -3531.4
$ws.Range("H84").Value = 2380.1667
$ws.Range("I84").Value = 2296.2
$ws.Range("K84").Value = 22962
$ws.Range("M84").Value = -17658
$ws.Range("H132").Value = 194162.03
$ws.Range("I132").Value = 1676.9574
$ws.Range("K132").Value = 5030.8722
$ws.Range("M132").Value = -2500.8722
